$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.043.04'
$ws.Range("E2").Value = '  +1.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.893.97'
$ws.Range("E3").Value = '  +1.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +1.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.30'
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4758'
$ws.Range("E7").Value = '  +1.72%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3972'
$ws.Range("E8").Value = '  +0.98%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.63'
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08054'
$ws.Range("E10").Value = '  +0.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.025'
$ws.Range("E11").Value = '  +0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.07'
$ws.Range("E12").Value = '  +1.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.889.60'
$ws.Range("E13").Value = '  +1.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.039'
$ws.Range("E14").Value = '  +1.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.245'
$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.016'
$ws.Range("E16").Value = '  +0.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.82'
$ws.Range("E17").Value = '  +2.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06788'
$ws.Range("E18").Value = '  +2.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001054'
$ws.Range("E19").Value = '  +0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.16'
$ws.Range("E20").Value = '  +0.18%  '

$ws.Range("E21").Value = '  +0.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.037.39'
$ws.Range("E22").Value = '  +1.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.557'
$ws.Range("E23").Value = '  +1.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.07'
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.351'
$ws.Range("E25").Value = '  +1.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.116.50'
$ws.Range("E26").Value = '  +1.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.06'
$ws.Range("E27").Value = '  +1.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.14'
$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.122'
$ws.Range("E29").Value = '  +1.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.569'
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.10'
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9842'
$ws.Range("E32").Value = '  +1.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09615'
$ws.Range("E33").Value = '  +1.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.647'
$ws.Range("E34").Value = '  +1.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.378'
$ws.Range("E35").Value = '  +1.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.376'
$ws.Range("E36").Value = '  -4.81%  '

# Row 37: coin swap (VeChain <-> Hedera) with updated data
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06112'
$ws.Range("E37").Value = '  +0.79%  '

# Row 38: coin swap (VeChain <-> Hedera) with updated data
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02264'
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.208'
$ws.Range("E39").Value = '  -2.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.253'
$ws.Range("E40").Value = '  +1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.013'
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6006'
$ws.Range("E42").Value = '  +0.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1902'
$ws.Range("E43").Value = '  +0.58%  '

$ws.Range("E44").Value = '  +1.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.278'
$ws.Range("E45").Value = '  +1.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5694'
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.23'
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.941'
$ws.Range("E48").Value = '  +0.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.375'
$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.65'
$ws.Range("E51").Value = '  -1.04%  '
